# Add field validator for verdict normalization in Truth model
# Updates verdict (column B) and confidence (column C) values for several rows
# to reflect the normalized verdict values and their associated confidences.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: confidence changes from 1 to 0.9 (verdict TRUE unchanged)
$ws.Range("C4").Value = 0.9

# Row 5: confidence changes from 0.8 to 1 (verdict FALSE unchanged)
$ws.Range("C5").Value = 1

# Row 6: verdict normalized from FALSE to INSUFFICIENT INFO; confidence cleared
$ws.Range("B6").Value = "INSUFFICIENT INFO"
$ws.Range("C6").Value = ""

# Row 8: verdict normalized from INSUFFICIENT INFO to FALSE; confidence set
# (leading apostrophe forces text instead of a boolean, then style is reset
# so no stray quote-prefix formatting is left behind)
$ws.Range("B8").Value = "'FALSE"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = 0.9333333333333333

# Row 9: verdict normalized from TRUE to FALSE; confidence changes from 1 to 0
$ws.Range("B9").Value = "'FALSE"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = 0

# Row 10: verdict normalized from TRUE to INSUFFICIENT INFO; confidence cleared
$ws.Range("B10").Value = "INSUFFICIENT INFO"
$ws.Range("C10").Value = ""

# Row 11: confidence changes from 1 to 0.8 (verdict TRUE unchanged)
$ws.Range("C11").Value = 0.8
